$wb = $excel.ActiveWorkbook

# --- Duplicate the "data" sheet into a new "dati" sheet (placed right after it) ---
$dataWs = $wb.Worksheets.Item("data")
$dataWs.Copy($null, $dataWs) | Out-Null

$datiWs = $wb.Worksheets.Item($wb.Worksheets.Count)
$datiWs.Name = "dati"

# --- Update the view/selection state on each sheet ---
$datiWs.Activate() | Out-Null
$datiWs.Range("C19").Select() | Out-Null

$dataWs.Activate() | Out-Null
$dataWs.Range("D16").Select() | Out-Null
